$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B2').NumberFormat = '@'
$ws.Range('B2').Value = '中性细胞数'
$ws.Range('C2').NumberFormat = '@'
$ws.Range('C2').Value = '6.3'
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '4-10'
$ws.Range('B3').NumberFormat = '@'
$ws.Range('B3').Value = '中性细胞数'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.5-5.5'
$ws.Range('C4').NumberFormat = '@'
$ws.Range('C4').Value = '93.0'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '110-160'
$ws.Range('C5').NumberFormat = '@'
$ws.Range('C5').Value = '31.3'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '36-50'
$ws.Range('C6').NumberFormat = '@'
$ws.Range('C6').Value = '67.5'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '86-100'
$ws.Range('B7').NumberFormat = '@'
$ws.Range('B7').Value = 'RDW-CV'
$ws.Range("C7").ClearContents() | Out-Null
$ws.Range("D7").ClearContents() | Out-Null
$ws.Range('B8').NumberFormat = '@'
$ws.Range('B8').Value = 'RDW-CV'
$ws.Range("C8").ClearContents() | Out-Null
$ws.Range("D8").ClearContents() | Out-Null
$ws.Range('B9').NumberFormat = '@'
$ws.Range('B9').Value = '血小板计数'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '264.00100-300'
$ws.Range("A10").ClearContents() | Out-Null
$ws.Range('B10').NumberFormat = '@'
$ws.Range('B10').Value = 'RDW-SL红细胞体积分布宽度-SD'
$ws.Range('C10').NumberFormat = '@'
$ws.Range('C10').Value = '41.2'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '37-50'
$ws.Range("A11").ClearContents() | Out-Null
$ws.Range('B11').NumberFormat = '@'
$ws.Range('B11').Value = 'RDW-CV'
$ws.Range('C11').NumberFormat = '@'
$ws.Range('C11').Value = '红细胞体积分布宽度-CV17.30↑'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '12-14.3'
$ws.Range('B12').NumberFormat = '@'
$ws.Range('B12').Value = '血小板分布宽度'
$ws.Range('C12').NumberFormat = '@'
$ws.Range('C12').Value = '13.0'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '9-17'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = '10.8'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9-13'
$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = '大型血小板比率'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value = '31.0'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '13-43'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.17-0.39'
$ws.Range('B16').NumberFormat = '@'
$ws.Range('B16').Value = 'RDW-CV'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2-7.7'
$ws.Range('A18').NumberFormat = '@'
$ws.Range('A18').Value = 'LYMPH#'
$ws.Range('B18').NumberFormat = '@'
$ws.Range('B18').Value = '淋巴细胞绝对值'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.8-4'
$ws.Range('C19').NumberFormat = '@'
$ws.Range('C19').Value = '1.11'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.25-0.95'
$ws.Range('B20').NumberFormat = '@'
$ws.Range('B20').Value = '嗜酸性粒细胞绝对值'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.01-0.59'
$ws.Range('B21').NumberFormat = '@'
$ws.Range('B21').Value = '嗜碱性粒细胞绝对值'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.01-0.07'
$ws.Range('B22').NumberFormat = '@'
$ws.Range('B22').Value = '中性粒细胞百分率'
$ws.Range('C22').NumberFormat = '@'
$ws.Range('C22').Value = '47.0'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '44-75'
$ws.Range('A23').NumberFormat = '@'
$ws.Range('A23').Value = 'LYMPH%淋巴细胞百分'
$ws.Range('B23').NumberFormat = '@'
$ws.Range('B23').Value = '淋巴细胞百分率'
$ws.Range('C23').NumberFormat = '@'
$ws.Range('C23').Value = '34.1'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '20-40'
$ws.Range('B24').NumberFormat = '@'
$ws.Range('B24').Value = '单核细胞百分率'
$ws.Range('C24').NumberFormat = '@'
$ws.Range('C24').Value = '17.6'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.2-15.2'
$ws.Range('B25').NumberFormat = '@'
$ws.Range('B25').Value = '嗜酸性粒细胞百分率'
$ws.Range("C25").ClearContents() | Out-Null
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.2-7.6'
$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = '嗜碱性粒细胞百分率'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1-1.2'
$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'RDW-CV'

# Cells that were already empty before the edit and remain empty (untouched by the diff);
# explicitly normalize them so the COM round-trip does not leave stray empty-string artifacts.
$ws.Range("C9").Value = $null
$ws.Range("A16").Value = $null
$ws.Range("C16").Value = $null
$ws.Range("C20").Value = $null
$ws.Range("C21").Value = $null
$ws.Range("C26").Value = $null
$ws.Range("A27").Value = $null
$ws.Range("C27").Value = $null
